$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.053.68"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "1.908.56"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'333.51"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.4638"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("D8").Value = "'0.4085"
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("D9").Value = "'47.95"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").Value = "'0.08029"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'1.007"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "'21.82"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "1.905.06"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "'5.956"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "'7.102"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'89.05"
$ws.Range("E16").Value = "  -2.28%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "'0.00001033"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").Value = "'0.06581"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "29.085.83"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").Value = "'5.452"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").Value = "2.132.91"
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("D27").Value = "'157.69"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").Value = "'19.74"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'2.114"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'5.417"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").Value = "'119.05"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").Value = "'0.9861"
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("D33").Value = "'0.09432"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").Value = "'1.425"
$ws.Range("E34").Value = "  +4.03%  "
$ws.Range("D35").Value = "'3.591"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").Value = "'5.326"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "'0.06098"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").Value = "'8.393"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "'1.176"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'0.5827"
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").Value = "'0.9994"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "'10.22"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").Value = "'0.1826"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("D46").Value = "'2.297"
$ws.Range("E46").Value = "  +11.11%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.18"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.5514"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'49.56"
$ws.Range("E49").Value = "  +27.74%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.917"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").Value = "'0.07045"
$ws.Range("E51").Value = "  +2.25%  "
